$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '62.872.27'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  +4.01%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.462.41'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  +3.64%  '
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  +0.18%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '408.86'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  -0.75%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '131.51'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  +18.47%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '3.454.24'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  +3.50%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.601'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  +2.51%  '
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  +0.15%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.692'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  +8.87%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.129'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  +29.25%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '42.94'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  +8.81%  '
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  -0.78%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '4.014.54'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  +3.86%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '8.78'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  +4.26%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '20.10'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  +1.93%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '3.449.87'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  +3.94%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '62.935.45'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  +4.72%  '
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  +0.08%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '10.82'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  +0.80%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.0000139'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  +26.33%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '3.33'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  -1.19%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '82.59'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  +10.33%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '13.14'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  +0.36%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '314.28'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  +4.28%  '
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  -0.69%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '30.37'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  +6.05%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '8.22'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  +2.46%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.61'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  -3.67%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '4.37'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  -2.44%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.177'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  -1.31%  '
$ws.Range("B32").NumberFormat = "@"
$ws.Range("B32").Value = 'InjectiveProtocol'
$ws.Range("C32").NumberFormat = "@"
$ws.Range("C32").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '44.37'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  +12.24%  '
$ws.Range("B33").NumberFormat = "@"
$ws.Range("B33").Value = 'Hedera'
$ws.Range("C33").NumberFormat = "@"
$ws.Range("C33").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.118'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  +2.25%  '
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  +2.53%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.60'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  -2.27%  '
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  +0.00%  '
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  -2.76%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '52.57'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  +0.38%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '3.57'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  +5.40%  '
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  +0.13%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '3.02'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  -3.10%  '
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  +1.87%  '
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  +3.82%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '135.74'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  -1.66%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '17.38'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  +2.69%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.286'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  -3.18%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.96'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  +0.58%  '
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  -0.13%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '22.07'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  -1.83%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '3.814.28'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  +3.93%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.180.55'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  +0.00%  '
